$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation")

$ws.Cells.Item(2, 2).Value = [double]"7.935379091409406E-19"
$ws.Cells.Item(2, 3).Value = [double]"1.8730848893942073E-65"
$ws.Cells.Item(2, 4).Value = [double]"2.1975287629784377E-67"

$ws.Cells.Item(3, 2).Value = [double]"2.9268971678437556E-5"
$ws.Cells.Item(3, 3).Value = [double]"11.850800084872"
$ws.Cells.Item(3, 4).Value = [double]"17.771731450919308"

$ws.Cells.Item(4, 2).Value = [double]"3.285787763735164E-5"
$ws.Cells.Item(4, 3).Value = [double]"24.016254020771903"
$ws.Cells.Item(4, 4).Value = [double]"49.784385219467374"

$ws.Cells.Item(5, 2).Value = [double]"3.109868213325136E-5"
$ws.Cells.Item(5, 3).Value = [double]"22.909527041942965"
$ws.Cells.Item(5, 4).Value = [double]"54.39191961738205"

$ws.Cells.Item(6, 2).Value = [double]"3.023658109532964E-5"
$ws.Cells.Item(6, 3).Value = [double]"21.53706482989181"
$ws.Cells.Item(6, 4).Value = [double]"57.41955020679926"

$ws.Cells.Item(7, 2).Value = [double]"3.0018253093708273E-5"
$ws.Cells.Item(7, 3).Value = [double]"18.544931734602798"
$ws.Cells.Item(7, 4).Value = [double]"61.30703126636506"

$ws.Cells.Item(8, 2).Value = [double]"3.0835954666303874E-5"
$ws.Cells.Item(8, 3).Value = [double]"13.141956925234286"
$ws.Cells.Item(8, 4).Value = [double]"66.70127869164776"

$ws.Cells.Item(9, 2).Value = [double]"3.201244665580381E-5"
$ws.Cells.Item(9, 3).Value = [double]"11.624496289983956"
$ws.Cells.Item(9, 4).Value = [double]"68.58750429666769"

$ws.Cells.Item(10, 2).Value = [double]"3.343245617309693E-5"
$ws.Cells.Item(10, 3).Value = [double]"11.266199862362518"
$ws.Cells.Item(10, 4).Value = [double]"69.0199062321509"

$ws.Cells.Item(11, 2).Value = [double]"4.001218782721183E-5"
$ws.Cells.Item(11, 3).Value = [double]"11.185668683451231"
$ws.Cells.Item(11, 4).Value = [double]"65.47358380173422"
